$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All values must be written as literal
# text (matching the original inline-string cells), even when they look numeric
# (e.g. '1.002', '30.774.55'), so Excel must not reinterpret them as numbers/dates.
$updates = [ordered]@{
    'D2' = '30.774.55'
    'E2' = '  -1.68%  '
    'D3' = '1.939.44'
    'E3' = '  -1.49%  '
    'D4' = '1.002'
    'E4' = '  +0.17%  '
    'D5' = '242.32'
    'E5' = '  -2.45%  '
    'E6' = '  +0.22%  '
    'D7' = '0.4884'
    'E7' = '  -0.32%  '
    'D8' = '0.2929'
    'E8' = '  -2.17%  '
    'D9' = '0.06881'
    'E9' = '  +0.14%  '
    'D10' = '19.52'
    'E10' = '  +0.94%  '
    'D11' = '105.42'
    'E11' = '  -1.86%  '
    'D12' = '1.956.09'
    'E12' = '  -0.50%  '
    'E13' = '  -0.40%  '
    'D14' = '5.314'
    'E14' = '  -3.09%  '
    'D15' = '0.6972'
    'E15' = '  -3.01%  '
    'D16' = '275.39'
    'E16' = '  -4.70%  '
    'D17' = '30.780.04'
    'E17' = '  -1.75%  '
    'B18' = 'Avalanche'
    'C18' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D18' = '13.15'
    'E18' = '  -1.55%  '
    'B19' = 'ShibaInu'
    'C19' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D19' = '0.000007707'
    'E19' = '  -1.03%  '
    'B20' = 'Dai'
    'C20' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D20' = '1.002'
    'E20' = '  +0.01%  '
    'B21' = 'WrappedliquidstakedEther2.0'
    'C21' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D21' = '2.194.95'
    'E21' = '  -0.54%  '
    'D22' = '5.437'
    'E22' = '  -3.94%  '
    'D23' = '1.001'
    'E23' = '  +0.14%  '
    'D24' = '6.460'
    'E24' = '  -3.12%  '
    'D25' = '9.704'
    'E25' = '  -3.70%  '
    'D26' = '167.59'
    'E26' = '  -1.15%  '
    'D27' = '19.57'
    'E27' = '  -2.52%  '
    'D28' = '2.161'
    'E28' = '  -1.89%  '
    'D29' = '0.1040'
    'E29' = '  -2.90%  '
    'D30' = '1.391'
    'E30' = '  -4.15%  '
    'D31' = '1.552'
    'E31' = '  -3.01%  '
    'D32' = '4.542'
    'E32' = '  -6.94%  '
    'D33' = '4.358'
    'E33' = '  -4.01%  '
    'D34' = '0.04846'
    'E34' = '  -4.84%  '
    'D35' = '0.7485'
    'E35' = '  -3.35%  '
    'D36' = '1.154'
    'E36' = '  -2.03%  '
    'D37' = '2.727'
    'E37' = '  -0.23%  '
    'D38' = '0.01987'
    'E38' = '  -3.65%  '
    'D39' = '2.663'
    'E39' = '  -1.99%  '
    'D40' = '77.36'
    'E40' = '  +4.67%  '
    'E41' = '  +0.05%  '
    'D42' = '2.091'
    'E42' = '  -3.15%  '
    'D43' = '0.8974'
    'E43' = '  +1.09%  '
    'D44' = '108.14'
    'E44' = '  -1.80%  '
    'D45' = '0.4406'
    'E45' = '  -2.17%  '
    'E46' = '  -0.19%  '
    'D47' = '7.740'
    'E47' = '  +2.38%  '
    'D48' = '989.00'
    'E48' = '  -0.75%  '
    'D49' = '0.1239'
    'E49' = '  -2.82%  '
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D50' = '35.73'
    'E50' = '  -1.10%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D51' = '9.155'
    'E51' = '  -2.80%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    $cell.Style = 'Normal'
}
